$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 students who dropped the course (entire-row delete so the
# remaining rows shift up and keep their relative order).
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Fill in the "Домашно 1" (column F) grades for week 6 for every remaining
# student.
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 5.4
$ws.Range("F4").Value = 5.6
$ws.Range("F5").Value = 5.4
$ws.Range("F6").Value = 3.75
$ws.Range("F7").Value = 5.5
$ws.Range("F8").Value = 5.25
$ws.Range("F9").Value = 3.5
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 5.95
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 6
$ws.Range("F14").Value = 5.95
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 4.95
$ws.Range("F18").Value = 6
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 5.4
$ws.Range("F21").Value = 6
$ws.Range("F22").Value = 5.7
$ws.Range("F23").Value = 3.5

# Match the author's final cursor position.
[void]$ws.Range("H16").Select()
